$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (headers): a new header value (6) is inserted at G1, shifting the old
# G1 value (7) into H1 and the old H1 value (10) into the new I1 cell.
$oldG1 = $ws.Range("G1").Value()
$oldH1 = $ws.Range("H1").Value()

$ws.Range("I1").Value = $oldH1
$ws.Range("H1").Value = $oldG1
$ws.Range("G1").Value = 6

# Copy H1's formatting (bold, centered, bordered header style) onto the new I1 cell
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Rows 2-10: add a new column I duplicating column H's value for that row
for ($r = 2; $r -le 10; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = $hVal
}
